$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "stream:datastream"
$ws.Range("B2").Value = "dict"

$ws.Range("A4").Value = "org:resource"
$ws.Range("B4").Value = "str"

$ws.Range("A5").Value = "concept:name"
$ws.Range("B5").Value = "str"

$ws.Range("A6").Value = "operation_end_time"
$ws.Range("B6").Value = "datetime"

$ws.Range("A7").Value = "SubProcessID"
$ws.Range("B7").Value = "str"
